$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 109, shifting existing data (rows 109-183) down to (111-185)
$ws.Rows("109:110").Insert()

# Populate new row 109
$ws.Range("A109").Value = 7
$ws.Range("B109").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C109").Value = "Ñuble"
$ws.Range("D109").Value = 44784
$ws.Range("E109").Value = 16
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100102
$ws.Range("H109").Value = "Cítricos"
$ws.Range("I109").Value = 100102004
$ws.Range("J109").Value = "Mandarina"
$ws.Range("K109").Value = "Murcott"
$ws.Range("L109").Value = "Primera"
$ws.Range("M109").Value = 60
$ws.Range("N109").Value = 8000
$ws.Range("O109").Value = 8000
$ws.Range("P109").Value = 8000
$ws.Range("Q109").Value = "`$/caja 18 kilos"
$ws.Range("R109").Value = "Región de O'Higgins"
$ws.Range("S109").Value = 444
$ws.Range("T109").Value = 18

# Populate new row 110
$ws.Range("A110").Value = 7
$ws.Range("B110").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C110").Value = "Ñuble"
$ws.Range("D110").Value = 44784
$ws.Range("E110").Value = 16
$ws.Range("F110").Value = "Fruta"
$ws.Range("G110").Value = 100102
$ws.Range("H110").Value = "Cítricos"
$ws.Range("I110").Value = 100102004
$ws.Range("J110").Value = "Mandarina"
$ws.Range("K110").Value = "Murcott"
$ws.Range("L110").Value = "Segunda"
$ws.Range("M110").Value = 120
$ws.Range("N110").Value = 7000
$ws.Range("O110").Value = 7500
$ws.Range("P110").Value = 7250
$ws.Range("Q110").Value = "`$/caja 18 kilos"
$ws.Range("R110").Value = "Región de O'Higgins"
$ws.Range("S110").Value = 403
$ws.Range("T110").Value = 18

# Apply date number format matching column D (style index 2: YYYY-MM-DD HH:MM:SS) to new D cells
$ws.Range("D109").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D110").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "Done"
